## update hydrogenation and improved dehydration conv ranges
##
## 1) Add a new validation column Q (rows 2-53) that flags whether each
##    parameter's "lower"/"upper" bound matches the expected relationship to
##    the baseline ($F$2), i.e. IF(E=H, 1, IF(F=$F$2, 1, 0)).
## 2) Replace row 37's (improved dehydration conversion) upper-bound formula
##    in column G with a symmetric-spread formula based on the new upper
##    bound in column I: 1-2*(I37-E37).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) New column Q: sanity-check flag for every parameter row ---------

# Row 2 is entered on its own (matches how the author typed the very first
# formula by hand before selecting the remainder of the column).
$ws.Range("Q2").Formula = "=IF(E2=H2, 1, IF(F2=`$F`$2, 1, 0))"

# Rows 3:53 are entered together so the engine groups them as one shared
# formula (mirrors dragging the fill handle down the rest of the column).
$ws.Range("Q3:Q53").Formula = "=IF(E3=H3, 1, IF(F3=`$F`$2, 1, 0))"

# --- 2) Row 37 (improved dehydration conversion): new upper-range formula -

$ws.Range("G37").Formula = "=1-2*(I37-E37)"

# --- 3) Selection/view bookkeeping matching the saved workbook state ------

# Try to scroll the window so row 18 is the first visible row (matches
# topLeftCell="A18" in the saved file); harmless if unsupported.
try {
    $excel.ActiveWindow.ScrollRow = 18
} catch {
}

# Select the entire row 37 (matches selection activeCell="A37"
# sqref="A37:XFD37").
$ws.Rows(37).Select()
